$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Tests")

# --- Sheet1: add new column R "Ex. Price" = ROUND(F,2) ---
$ws1.Range("R1").Value = "Ex. Price"
$ws1.Range("R2").Formula = "=ROUND(F2,2)"
$ws1.Range("R3:R49").Formula = "=ROUND(F3,2)"

# Update the selection on Sheet1 to the new column R
$ws1.Range("R1:R1048576").Select()

# --- Tests sheet: fix two rounded Ex. Price literal values ---
$ws2.Range("F4").Value = 58.35
$ws2.Range("F8").Value = 86.83

# Update the selection on Tests sheet
$ws2.Range("P10").Select()

# Re-activate the Tests sheet (it was the tab selected/active one)
$ws2.Activate()
